$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.933.22"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3
$ws.Range("D3").Value = "1.887.47"
$ws.Range("E3").Value = "  +0.91%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.42"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4570"
$ws.Range("E7").Value = "  +0.67%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3902"
$ws.Range("E8").Value = "  +1.47%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07856"
$ws.Range("E9").Value = "  +0.29%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9864"
$ws.Range("E10").Value = "  -0.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.91"
$ws.Range("E11").Value = "  +1.87%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.926.64"
$ws.Range("E12").Value = "  +2.87%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.036"
$ws.Range("E13").Value = "  +1.67%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.694"
$ws.Range("E14").Value = "  +0.89%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06947"
$ws.Range("E15").Value = "  +0.12%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.05"
$ws.Range("E16").Value = "  +1.70%  "

# Row 17
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009980"
$ws.Range("E18").Value = "  +0.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.02"
$ws.Range("E19").Value = "  +1.80%  "

# Row 20
$ws.Range("E20").Value = "  -0.36%  "

# Row 21
$ws.Range("D21").Value = "28.920.69"
$ws.Range("E21").Value = "  +1.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.286"
$ws.Range("E22").Value = "  +0.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  +0.59%  "

# Row 24
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.060.08"
$ws.Range("E24").Value = "  -2.02%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.056"
$ws.Range("E25").Value = "  -0.65%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.73"
$ws.Range("E26").Value = "  +1.23%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.33"
$ws.Range("E27").Value = "  +1.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.942"
$ws.Range("E28").Value = "  +4.72%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.926"
$ws.Range("E29").Value = "  +2.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.59"
$ws.Range("E30").Value = "  +0.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09326"
$ws.Range("E31").Value = "  +0.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9072"
$ws.Range("E32").Value = "  +0.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.292"
$ws.Range("E33").Value = "  +0.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.330"
$ws.Range("E34").Value = "  +0.70%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.262"
$ws.Range("E35").Value = "  +0.16%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.203"
$ws.Range("E36").Value = "  +4.72%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05765"
$ws.Range("E37").Value = "  +1.73%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02069"
$ws.Range("E38").Value = "  +1.41%  "

# Row 39
$ws.Range("E39").Value = "  -0.26%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5675"
$ws.Range("E40").Value = "  +1.90%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.637"
$ws.Range("E41").Value = "  +0.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1770"
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.741"
$ws.Range("E43").Value = "  +0.89%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.281"
$ws.Range("E44").Value = "  +6.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.95"
$ws.Range("E45").Value = "  +4.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5363"
$ws.Range("E46").Value = "  +2.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07039"
$ws.Range("E47").Value = "  -1.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.854"
$ws.Range("E48").Value = "  +2.74%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.49"
$ws.Range("E49").Value = "  +0.72%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.526"
$ws.Range("E50").Value = "  +3.56%  "

# Row 51
$ws.Range("E51").Value = "  -3.58%  "
